# Add custom widgets (custom_date / custom_date_time) to the example
# survey template: introduce a "clause" column + a templatePath column
# on the survey sheet, wrap the two new widget rows in a begin/end
# screen group, and rename the form's display title in settings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# survey sheet: insert the "clause" / "templatePath" columns and the
# begin-screen / custom widgets / end-screen rows.
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item(1)

# Header row: shift the existing headers one column to the right and
# add the new "clause" (first) and "templatePath" (last) columns.
$survey.Cells.Item(1,1).Value = "clause"
$survey.Cells.Item(1,2).Value = "type"
$survey.Cells.Item(1,3).Value = "values_list"
$survey.Cells.Item(1,4).Value = "name"
$survey.Cells.Item(1,5).Value = "display.prompt.text"
$survey.Cells.Item(1,6).Value = "templatePath"
$survey.Rows.Item(1).RowHeight = 25.5

# Row 2 becomes the "begin screen" clause row (name/value cells cleared).
$survey.Cells.Item(2,1).Value = "begin screen"
$survey.Cells.Item(2,2).ClearContents()
$survey.Cells.Item(2,3).ClearContents()
$survey.Cells.Item(2,4).ClearContents()

# Row 3: custom_date widget backed by a handlebars template.
$survey.Cells.Item(3,2).Value = "text"
$survey.Cells.Item(3,4).Value = "custom_date"
$survey.Cells.Item(3,5).Value = "Specify Date"
$survey.Cells.Item(3,6).Value = "../../../_templates/custom_date_picker.handlebars"
$survey.Rows.Item(3).RowHeight = 38.25

# Row 4: custom_date_time widget backed by a handlebars template.
$survey.Cells.Item(4,2).Value = "text"
$survey.Cells.Item(4,4).Value = "custom_date_time"
$survey.Cells.Item(4,5).Value = "Specify Date & Time"
$survey.Cells.Item(4,6).Value = "../../../_templates/custom_datetime_picker.handlebars"
$survey.Rows.Item(4).RowHeight = 38.25

# Row 5: close the screen clause group.
$survey.Cells.Item(5,1).Value = "end screen"

# ---------------------------------------------------------------------
# settings sheet: rename the form's display title.
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item(3)
$settings.Cells.Item(5,3).Value = "Templates Example Form"

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping to match the authored workbook.
# ---------------------------------------------------------------------
$settings.Activate()
$settings.Range("C6").Select()

$survey.Activate()
$survey.Range("A11").Select()
